# Insert a new "MAE" column before the existing "Tipo" column and update
# the regression metric values for the AdaBoostRegressor prediction sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "Tipo" column (D) one to the right to make room for "MAE".
$ws.Columns.Item(4).Insert()

# New header for the inserted column, matching the style of the other headers.
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Updated metric values.
$ws.Range("B2").Value = 0.2405818208600997
$ws.Range("C2").Value = 0.9952978222611105
$ws.Range("D2").Value = 0.380904464225422
